$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Sheet2: brand new "Raster:" comparison table (Raster vs Buffer LatLongBox),
# built from scratch to mirror the target layout.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Column widths (best achievable approximation of the target "best fit"
# widths - 13.71, 21.71, 19.86, 22.71, 17.29, 14.71 characters).
$ws2.Columns.Item(1).ColumnWidth = 12.9
$ws2.Columns.Item(2).ColumnWidth = 20.9
$ws2.Columns.Item(3).ColumnWidth = 19.05
$ws2.Columns.Item(4).ColumnWidth = 21.9
$ws2.Columns.Item(5).ColumnWidth = 16.45
$ws2.Columns.Item(6).ColumnWidth = 13.9

# Header row + row labels (written in the same order the strings are first
# introduced so the shared-string table lines up with the target layout).
$ws2.Range("C1").Value = "Raster:"
$ws2.Range("A2").Value = "North"
$ws2.Range("A3").Value = "South"
$ws2.Range("A4").Value = "East"
$ws2.Range("A5").Value = "West"
$ws2.Range("D1").Value = "Raster LatLongBox"
$ws2.Range("E1").Value = "Buffer LatLongBox"
$ws2.Range("A10").Value = "Tile Lat Height"
$ws2.Range("A11").Value = "Tile Lon Width"
$ws2.Range("B1").Value = "Model"

# Numeric data: column B keeps the 0.00000000000000000-style (17-dp) format
# used elsewhere in the workbook; columns C/D/E use the new finer
# 0.000000000000 (12-dp) format.
$ws2.Range("B2:B5").NumberFormat = "0.00000000000000000"
$ws2.Range("C2:E5").NumberFormat = "0.000000000000"
$ws2.Range("C10:C11").NumberFormat = "0.000000000000"
$ws2.Range("B7:B8").NumberFormat = "0.00000000000000000"
$ws2.Range("C7:C8").NumberFormat = "0.000000000000"

$ws2.Range("B2").Value = 44.323888888653777
$ws2.Range("C2").Value = 44.285831541747001
$ws2.Range("D2").Value = 44.285831541747001
$ws2.Range("E2").Value = 44.231388332300398

$ws2.Range("B3").Value = 44.211666666424001
$ws2.Range("C3").Value = 44.239627838134702
$ws2.Range("D3").Value = 44.239627838134702
$ws2.Range("E3").Value = 44.211758699770101

$ws2.Range("B4").Value = -71.244722218675577
$ws2.Range("C4").Value = -71.279353660393099
$ws2.Range("D4").Value = -71.279353660393099
$ws2.Range("E4").Value = -71.292219848815805

$ws2.Range("B5").Value = -71.384722218684999
$ws2.Range("C5").Value = -71.325927734375
$ws2.Range("D5").Value = -71.325927734375
$ws2.Range("E5").Value = -71.384719848632798

$ws2.Range("B7").Formula = "=B2-B3"
$ws2.Range("C7").Formula = "=C2-C3"

$ws2.Range("B8").Formula = "=B4-B5"
$ws2.Range("C8").Formula = "=C4-C5"

$ws2.Range("C10").Value = 0.092499999816936906
$ws2.Range("C11").Value = 0.092499999816936906

$ws2.Range("B4").Select() | Out-Null

# ---------------------------------------------------------------------------
# Restore Sheet1 as the active/visible tab and move its selection from D9 to
# E9, matching the target workbook view state.
# ---------------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("E9").Select() | Out-Null

$wb.Save()
